$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: <one> -> <india>, count 9 -> 7
$ws.Range("A2").Value = "<india>"
$ws.Range("B2").Value = "<india>"
$ws.Range("C2").Value = 7

# Row 3: <oil> -> <uniform>, count 5 (unchanged)
$ws.Range("A3").Value = "<uniform>"
$ws.Range("B3").Value = "<uniform>"

# Row 4: <up> -> <oscar>, count 10 -> 6
$ws.Range("A4").Value = "<oscar>"
$ws.Range("B4").Value = "<oscar>"
$ws.Range("C4").Value = 6

# Row 5: <up> -> <water>, count 17 -> 4
$ws.Range("A5").Value = "<water>"
$ws.Range("B5").Value = "<water>"
$ws.Range("C5").Value = 4

# Row 6: <like> -> <so>, count 10 -> 7
$ws.Range("A6").Value = "<so>"
$ws.Range("B6").Value = "<so>"
$ws.Range("C6").Value = 7

# Row 7: A unchanged <and>, B <an> -> <and>, count 14 -> 7
$ws.Range("B7").Value = "<and>"
$ws.Range("C7").Value = 7

# Row 8: <shift> -> <zero>, count 12 -> 6
$ws.Range("A8").Value = "<zero>"
$ws.Range("B8").Value = "<zero>"
$ws.Range("C8").Value = 6

# Row 9: <said> -> <when>, count 13 -> 4
$ws.Range("A9").Value = "<when>"
$ws.Range("B9").Value = "<when>"
$ws.Range("C9").Value = 4

# Row 10: A <tango> -> <sentence>, B <hin> -> <sentence>, count 12 -> 4
$ws.Range("A10").Value = "<sentence>"
$ws.Range("B10").Value = "<sentence>"
$ws.Range("C10").Value = 4

# Row 11: <mike> -> <cut>, count 14 -> 7
$ws.Range("A11").Value = "<cut>"
$ws.Range("B11").Value = "<cut>"
$ws.Range("C11").Value = 7

# Row 12: <kilo> -> <could>, count 8 -> 5
$ws.Range("A12").Value = "<could>"
$ws.Range("B12").Value = "<could>"
$ws.Range("C12").Value = 5

# Row 13: <a> -> <delta>, count 10 (unchanged)
$ws.Range("A13").Value = "<delta>"
$ws.Range("B13").Value = "<delta>"

# Row 14: <may> -> <that>, count 9 -> 6
$ws.Range("A14").Value = "<that>"
$ws.Range("B14").Value = "<that>"
$ws.Range("C14").Value = 6

# Row 15: <water> -> <first>, count 12 -> 4
$ws.Range("A15").Value = "<first>"
$ws.Range("B15").Value = "<first>"
$ws.Range("C15").Value = 4

# Row 16: <zulu> -> <could>, count 9 -> 4
$ws.Range("A16").Value = "<could>"
$ws.Range("B16").Value = "<could>"
$ws.Range("C16").Value = 4

# Row 17: <water> -> <then>, count 11 -> 5
$ws.Range("A17").Value = "<then>"
$ws.Range("B17").Value = "<then>"
$ws.Range("C17").Value = 5

# Row 18: A <uniform> -> <can>, B <in> -> <can>, count 15 -> 6
$ws.Range("A18").Value = "<can>"
$ws.Range("B18").Value = "<can>"
$ws.Range("C18").Value = 6
